$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 25
    $ws.Range("F4").Value = 6204
    $ws.Range("F8").Value = 1867
    $ws.Range("F12").Value = 215
}
